$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (losing trailing zeros / thousands-dot formatting) are forced to Text format
# before the value is written, so the literal string is preserved verbatim.

$ws.Range("D2").Value = "61.944.26"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "3.398.47"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.74"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.50"
$ws.Range("E6").Value = "  +3.16%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.396.55"
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.127"
$ws.Range("E11").Value = "  +3.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.392"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "3.983.80"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").Value = "3.395.58"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.47"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "62.045.28"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.19"
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.51"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "394.49"
$ws.Range("E22").Value = "  +3.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.567"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000130"
$ws.Range("E24").Value = "  +9.37%  "
$ws.Range("D25").Value = "3.547.21"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.66"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.66"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.161"
$ws.Range("E31").Value = "  +3.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.24"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.52"
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("D36").Value = "3.434.83"
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.42"
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("E38").Value = "  +2.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.92"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "164.93"
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0790"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.78"
$ws.Range("E42").Value = "  +9.82%  "
$ws.Range("B43").Value = "ONDO"
$ws.Range("C43").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.26"
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.790"
$ws.Range("E44").Value = "  +4.07%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.44"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.07"
$ws.Range("E47").Value = "  +5.93%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.39"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.91"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.19"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "2.342.71"
$ws.Range("E51").Value = "  +6.75%  "
